$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")
$ws.Activate()

$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("O2").ClearContents()

$ws.Range("M2:O2").Select()
